$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25, column A ("71076783") was stored as text; the redemption is now
# recorded as a numeric phone id.
$ws.Range("A25").Value = 71076783

# Append the new redemption event as row 26: phone (text, matches existing
# inlineStr convention for the phone column), points, timestamp.
$ws.Range("A26").NumberFormat = "@"
$ws.Range("A26").Value = "71076783"
$ws.Range("A26").ClearFormats()

$ws.Range("B26").Value = 100
$ws.Range("C26").Value = "2025-08-18T18:01:45"
